$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tags")

$ws.Range("B21").Value = "atomizer-type"
$ws.Range("A21").Value = "RTA"
$ws.Range("C21").Value = 0

$ws.Range("B22").Value = "atomizer-type"
$ws.Range("A22").Value = "RDA"
$ws.Range("C22").Value = 1

$ws.Range("B23").Value = "atomizer-type"
$ws.Range("A23").Value = "RDTA"
$ws.Range("C23").Value = 2

$ws.Range("A23").Select()
